$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "City" value for the Chennai row to "Delhi"
$ws.Range("A2").Value = "Delhi"

# Add a new row with "London" in A3
$ws.Range("A3").Value = "London"

# Move the active selection to B4 (single cell, not a range)
$ws.Range("B4").Select()
